# Bloch sphere slide: add the |0> / |1> state labels near the poles of the
# sphere, and tidy the stray endParaRPr on the "y" axis label text box.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# --- locate existing shapes by name (robust to index drift) -----------------
$yShape = $null
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $sh = $s.Shapes.Item($i)
    if ($sh.Name -eq "Textfeld 55") { $yShape = $sh }
}

# Drop the redundant trailing endParaRPr on the "y" label's paragraph: clear
# the text (this also clears the stray endParaRPr) and retype it, which keeps
# the run's existing formatting (italic Times New Roman) intact.
if ($yShape -ne $null) {
    $yShape.TextFrame.DeleteText()
    $yShape.TextFrame.TextRange.Text = "y"
}

# --- add the |0> label near the top pole ------------------------------------
$ket0 = $s.Shapes.AddTextbox(1, 3635896/12700, 1484784/12700, 522900/12700, 369332/12700)
$ket0.Name = "Textfeld 27"
$ket0.TextFrame.WordWrap = 0
$ket0.TextFrame.AutoSize = 1
$ket0.Fill.Visible = 0

$tr0 = $ket0.TextFrame.TextRange
$tr0.Text = "|0>"
$tr0.LanguageID = "de-DE"
$tr0.Font.Name = "Arial"
$tr0.Font.NameComplexScript = "Arial"

# --- add the |1> label near the bottom pole ---------------------------------
$ket1 = $s.Shapes.AddTextbox(1, 3617052/12700, 4293096/12700, 522900/12700, 369332/12700)
$ket1.Name = "Textfeld 29"
$ket1.TextFrame.WordWrap = 0
$ket1.TextFrame.AutoSize = 1
$ket1.Fill.Visible = 0

$tr1 = $ket1.TextFrame.TextRange
$tr1.Text = "|1>"
$tr1.LanguageID = "de-DE"
$tr1.Font.Name = "Arial"
$tr1.Font.NameComplexScript = "Arial"
